# Update "Forecast Comparison" sheet with a new Week_Start_Date column,
# shortened Week labels, a boolean is_holiday_week column, and corrected
# MyForecast values for W10 and W16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current ASIN column (B) to hold
# "Week_Start_Date". This shifts ASIN..is_holiday_week one column right.
$ws.Columns.Item(2).Insert()

# --- Header row ---
$ws.Range("B1").Value = "Week_Start_Date"

# --- Week labels (A2:A17): drop the leading zero, e.g. W01 -> W1 ---
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $weeks[$i]
}

# --- Week_Start_Date values (B2:B17), kept as plain text ---
$ws.Range("B2:B17").NumberFormat = "@"
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)
for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
}

# --- Corrected MyForecast values (column D after the insert) ---
# Row 11 (W10): 4 -> 5
$ws.Cells.Item(11, 4).Value = 5
# Row 17 (W16): 5 -> 4
$ws.Cells.Item(17, 4).Value = 4

# --- is_holiday_week (column J after the insert) becomes boolean ---
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 10).Value = $false
}
